$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: extend the stair-climbing instruction text.
$ws.Range("C3").Value = '"Suba las escaleras principales hasta el  segundo descanso, llegará al segundo piso."'

# Row 4: "concejo" was a misspelling of "consejo".
$ws.Range("B4").Value = "consejo"

# Row 3: shorten the warning text.
$ws.Range("D3").Value = '"escalones resbaladizos"'

# The row now wraps onto more lines, so its height grows (matches the autosized
# height already used by row 5, which holds similarly long wrapped text).
$ws.Range("A3").EntireRow.RowHeight = 40.2

# Row 7: add a new "informacion" destination summarizing the offices reachable.
$ws.Range("B7").Value = "informacion"
$ws.Range("C7").Value = " consejo, despacho, cobro coactivo"
